$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header-style rows below the existing data, matching the
# formatting used by the other single-column header rows (e.g. A1).
$ws.Range("A1").Copy($ws.Range("A25"))
$ws.Range("A25").Value = "DfT Group"

$ws.Range("A1").Copy($ws.Range("A26"))
$ws.Range("A26").Value = "IPDC approval point"

# Move/update the active selection to A27
$ws.Range("A27").Select()
